$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Row 4: new script entry "Login" - Script #1
$ws.Range("A4").Value = "1"
$ws.Range("B4").Value = "Login"
$ws.Range("C4").Value = "2024-04-30"
$ws.Range("D4").Value = "EN"
$ws.Range("E4").Value = "PASS"

# Row 11: previously held "8 / Load or Pay credit card / ..." - now cleared (content only, keep formatting)
$ws.Range("A11:E11").ClearContents()

# Row 16: previously held "13 / Check net worth statment chasrts / ..." - now cleared (content only, keep formatting)
$ws.Range("A16:E16").ClearContents()

# Row 17: new script entry "Check wazin account details" - Script #14 - FAIL
$ws.Range("A17").Value = "14"
$ws.Range("B17").Value = "Check wazin account details "
$ws.Range("C17").Value = "2024-04-30"
$ws.Range("D17").Value = "EN"
$ws.Range("E17").Value = " FAIL "

# Row 18: new script entry "Between my account transfer" - Script #15 - PASS
$ws.Range("A18").Value = "15"
$ws.Range("B18").Value = "Between my account transfer"
$ws.Range("C18").Value = "2024-04-30"
$ws.Range("D18").Value = "EN"
$ws.Range("E18").Value = "PASS"

# Row 19: new script entry "Within riyad bank trnasfer" - Script #16 - PASS
$ws.Range("A19").Value = "16"
$ws.Range("B19").Value = "Within riyad bank trnasfer"
$ws.Range("C19").Value = "2024-04-30"
$ws.Range("D19").Value = "EN"
$ws.Range("E19").Value = "PASS"

# Update the active cell selection
$ws.Range("E12").Select()

$wb.Save()
